# Insert a new data row before row 63 (the existing rows 63..120 shift down
# to 64..121). The new row duplicates the values of what was row 98
# ("2023-02-07" / serial 44964), except the date, which becomes serial 44966
# ("2023-02-09").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(63).Insert()

$newRow = 63

$ws.Cells.Item($newRow, 1).Value = 7
$ws.Cells.Item($newRow, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($newRow, 3).Value = "Ñuble"
$ws.Cells.Item($newRow, 4).Value = 44966
$ws.Cells.Item($newRow, 5).Value = 16
$ws.Cells.Item($newRow, 6).Value = 100112031
$ws.Cells.Item($newRow, 7).Value = "Poroto verde"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 60
$ws.Cells.Item($newRow, 11).Value = 26000
$ws.Cells.Item($newRow, 12).Value = 28000
$ws.Cells.Item($newRow, 13).Value = 27000
$ws.Cells.Item($newRow, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item($newRow, 15).Value = "Región del Maule"
$ws.Cells.Item($newRow, 16).Value = 1080
$ws.Cells.Item($newRow, 17).Value = 25
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
